$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The list was re-uploaded with one fewer row: the entries that used to sit
# in rows 11 and 12 move up into rows 10 and 11, and the old last row (12)
# goes away.
#
# Force the destination cells to text first so the phone numbers (which
# start with "+") and the ISO dates are stored as literal text instead of
# being auto-converted to a number/date by Excel, matching how the rest of
# the sheet stores its data.
$ws.Range("A10:C11").NumberFormat = "@"

$ws.Range("A10").Value = "+5511947261969"
$ws.Range("B10").Value = "11"
$ws.Range("C10").Value = "2024-09-10"

$ws.Range("A11").Value = "+5521981400589"
$ws.Range("B11").Value = "21"
$ws.Range("C11").Value = "2024-09-09"

# Re-apply the same formatting (font/alignment/border/General number format)
# used by every other data row so these two rows don't end up on a
# different, newly-minted cell style than the rest of the table.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C11").PasteSpecial(-4122)  # xlPasteFormats

# Row 12 (the old last entry, now duplicated into row 11 above) is removed.
$ws.Rows.Item(12).Delete()
